$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q3 at the
#    top of the data (row 2), shifting the existing quarters down, and
#    renumber the index column (A) to stay 0..n-1.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

$wsTotal.Rows.Item(2).Insert()
$wsTotal.Rows.Item(2).ClearFormats()

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 21
$wsTotal.Cells.Item(2, 4).Value = 4.33

for ($r = 3; $r -le 9; $r++) {
    $wsTotal.Cells.Item($r, 1).Value = $r - 2
}

# Restore the index-column style (bold/border) on the new row's A cell by
# copying it from a row that still carries it.
$wsTotal.Cells.Item(3, 1).Copy()
$wsTotal.Cells.Item(2, 1).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Insert a brand-new worksheet "2022-Q3" right after "总计" (i.e. before
#    the current "2022-Q2" tab) holding the per-fund holdings breakdown.
#    Duplicate the "2022-Q2" sheet (so sheetPr/pageMargins/etc. match the
#    rest of the workbook) then wipe its contents before refilling it.
# ---------------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item(2)
$wsQ2.Copy($wsQ2)
$wsNew = $wb.Worksheets.Item(2)
$wsNew.Name = "2022-Q3"
$wsNew.Cells.Clear()

# Header row (B1:H1), reusing the bold/border header style from "总计"!B1.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$wsTotal.Cells.Item(1, 2).Copy()
$wsNew.Range("B1:H1").PasteSpecial(-4122)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsNew.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$fundData = @(
    @("0","000727","融通健康产业灵活配置混合A","22.64","93.67","6.12","1.3856","4"),
    @("1","009274","融通健康产业灵活配置混合C","17.64","93.67","6.12","1.0796","4"),
    @("2","012159","财通资管健康产业混合A","10.00","94.52","5.50","0.5500","7"),
    @("3","012173","国泰兴泽优选一年持有期混合A","8.41","88.23","3.32","0.2792","8"),
    @("4","011466","兴业医疗保健混合A","3.85","87.15","5.91","0.2275","4"),
    @("5","012174","国泰兴泽优选一年持有期混合C","6.17","88.23","3.32","0.2048","8"),
    @("6","001551","天弘中证医药100指数型发起式 C","8.58","95.24","1.48","0.1270","4"),
    @("7","011467","兴业医疗保健混合C","1.94","87.15","5.91","0.1147","4"),
    @("8","012160","财通资管健康产业混合C","1.95","94.52","5.50","0.1072","7"),
    @("9","011404","融通鑫新成长混合C","1.75","94.07","4.71","0.0824","5"),
    @("10","001550","天弘中证医药100指数型发起式 A","5.31","95.24","1.48","0.0786","4"),
    @("11","013441","西藏东财创新医疗六个月定开混合","0.49","82.53","5.06","0.0248","6"),
    @("12","011403","融通鑫新成长混合A","0.39","94.07","4.71","0.0184","5"),
    @("13","008619","永赢医药健康股票C","0.40","94.40","4.39","0.0176","9"),
    @("14","008618","永赢医药健康股票A","0.24","94.40","4.39","0.0105","9"),
    @("15","014462","光大保德信汇佳混合A","0.33","43.38","2.76","0.0091","5"),
    @("16","013920","兴华创新医疗6个月持有混合A","0.18","94.83","4.07","0.0073","8"),
    @("17","005105","富荣福康混合C","0.07","91.00","3.23","0.0023","6"),
    @("18","013921","兴华创新医疗6个月持有混合C","0.05","94.83","4.07","0.0020","8"),
    @("19","005104","富荣福康混合A","0.03","91.00","3.23","0.0010","6"),
    @("20","014463","光大保德信汇佳混合C","0.03","43.38","2.76","0.0008","5"),
)

for ($i = 0; $i -lt $fundData.Length; $i++) {
    $r = $fundData[$i]
    $row = $i + 2

    $wsNew.Cells.Item($row, 1).Value = [int]$r[0]
    Set-TextValue $wsNew.Cells.Item($row, 2) $r[1]
    Set-TextValue $wsNew.Cells.Item($row, 3) $r[2]
    Set-TextValue $wsNew.Cells.Item($row, 4) $r[3]
    Set-TextValue $wsNew.Cells.Item($row, 5) $r[4]
    Set-TextValue $wsNew.Cells.Item($row, 6) $r[5]
    Set-TextValue $wsNew.Cells.Item($row, 7) $r[6]
    $wsNew.Cells.Item($row, 8).Value = [int]$r[7]
}

# Restore the index-column style (bold/border) on column A for the data rows.
$wsTotal.Cells.Item(2, 1).Copy()
$wsNew.Range("A2:A22").PasteSpecial(-4122)

# Keep the originally-active sheet/selection ("总计"!A1) selected, since the
# sheet-insert/copy operations above shift focus onto the new tab.
$wsTotal.Activate()
$wsTotal.Range("A1").Select() | Out-Null

Write-Output "edit complete"
